$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "standard deviation" rows (41-43) - labels first
$ws.Range("B41").Value = "Standard deviation of KP"
$ws.Range("D41").Formula = "=STDEV.S(D6:D25)"

$ws.Range("B42").Value = "Standard deviation of KPB"
$ws.Range("D42").Formula = "=STDEV.S(G6:G25)"

$ws.Range("B43").Value = "Standard deviation of KPD"
$ws.Range("D43").Formula = "=STDEV.S(J6:J25)"

# New "95% confidence interval" rows (45-47) - labels first
$ws.Range("B45").Value = "95% Confidence Interval KP"
$ws.Range("D45").Formula = "=CONFIDENCE.NORM(0.05, D41, 20)"

$ws.Range("B46").Value = "95% Confidence Interval KPB"
$ws.Range("D46").Formula = "=CONFIDENCE.NORM(0.05,D42,20)"

$ws.Range("B47").Value = "95% Confidence Interval KPD"
$ws.Range("D47").Formula = "=CONFIDENCE.NORM(0.05, D43,20)"

# Column-A "Synthetic" tag, added last to match the original authoring order
$ws.Range("A41").Value = "Synthetic"
$ws.Range("A42").Value = "Synthetic"
$ws.Range("A43").Value = "Synthetic"
$ws.Range("A45").Value = "Synthetic"
$ws.Range("A46").Value = "Synthetic"
$ws.Range("A47").Value = "Synthetic"

# Update the selected cell to match the post-edit cursor position
$ws.Range("A48").Select() | Out-Null
